$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4902.25
$ws.Range("I51").Value = 4438
$ws.Range("K51").Value = 4438
$ws.Range("M51").Value = -3954

$ws.Range("H98").Value = 2312.8438
$ws.Range("I98").Value = 2333.5667
$ws.Range("K98").Value = 2333.5667
$ws.Range("M98").Value = -835.5666999999999

$ws.Range("H112").Value = 4011.9285
$ws.Range("I112").Value = 966.3333
$ws.Range("K112").Value = 2898.9999
$ws.Range("M112").Value = -1790.9999

$ws.Range("H122").Value = 2312.8438
$ws.Range("I122").Value = 2333.5667
$ws.Range("K122").Value = 7000.7001
$ws.Range("M122").Value = -4550.7001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3035917.5
$ws.Range("I32").Value = 3081854.8
$ws.Range("J32").Value = 49999
$ws.Range("K32").Value = 3081854.8
$ws.Range("L32").Value = 49999
$ws.Range("M32").Value = -3081567.8
$ws.Range("N32").Value = -50573

$ws.Range("H45").Value = 4535.1313
$ws.Range("I45").Value = 3195.2
$ws.Range("J45").Value = 5409
$ws.Range("K45").Value = 3195.2
$ws.Range("L45").Value = 5409
$ws.Range("M45").Value = -2818.2
$ws.Range("N45").Value = -6163

$ws.Range("H98").Value = 47087.5
$ws.Range("J98").Value = 47087.5
$ws.Range("L98").Value = 47087.5
$ws.Range("N98").Value = -53077.5

$ws.Range("H102").Value = 3024.6956
$ws.Range("I102").Value = 2187.842
$ws.Range("K102").Value = 2187.842
$ws.Range("M102").Value = -565.8420000000001

$ws.Range("H110").Value = 37038436
$ws.Range("I110").Value = 1553.8572
$ws.Range("K110").Value = 1553.8572
$ws.Range("M110").Value = 491.1428000000001

$ws.Range("H122").Value = 3042.0952
$ws.Range("I122").Value = 2706.9412
$ws.Range("K122").Value = 8120.823600000001
$ws.Range("M122").Value = -5670.823600000001

$ws.Range("H131").Value = 51586.5
$ws.Range("J131").Value = 51586.5
$ws.Range("L131").Value = 51586.5
$ws.Range("N131").Value = -61666.5

$ws.Range("H132").Value = 3724.9768
$ws.Range("I132").Value = 1247.32
$ws.Range("J132").Value = 7166.1665
$ws.Range("K132").Value = 3741.96
$ws.Range("L132").Value = 21498.4995
$ws.Range("M132").Value = -1211.96
$ws.Range("N132").Value = -26558.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5323590
$ws.Range("I134").Value = 7578776
$ws.Range("K134").Value = 22736328
$ws.Range("M134").Value = -22733793

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5240.024
$ws.Range("I31").Value = 1933.0625
$ws.Range("J31").Value = 7275.077
$ws.Range("K31").Value = 1933.0625
$ws.Range("L31").Value = 7275.077
$ws.Range("M31").Value = -1638.0625
$ws.Range("N31").Value = -7865.077

$ws.Range("H34").Value = 5240.024
$ws.Range("I34").Value = 1933.0625
$ws.Range("J34").Value = 7275.077
$ws.Range("K34").Value = 1933.0625
$ws.Range("L34").Value = 7275.077
$ws.Range("M34").Value = -1731.0625
$ws.Range("N34").Value = -7679.077

$ws.Range("H99").Value = 4395.25
$ws.Range("I99").Value = 1820.4286
$ws.Range("K99").Value = 1820.4286
$ws.Range("M99").Value = -322.4286

$ws.Range("H126").Value = 4395.25
$ws.Range("I126").Value = 1820.4286
$ws.Range("K126").Value = 5461.2858
$ws.Range("M126").Value = -2991.2858

$ws.Range("H132").Value = 4821.9165
$ws.Range("I132").Value = 3579.1
$ws.Range("J132").Value = 6375.4375
$ws.Range("K132").Value = 10737.3
$ws.Range("L132").Value = 19126.3125
$ws.Range("M132").Value = -8207.299999999999
$ws.Range("N132").Value = -24186.3125

$ws.Range("H134").Value = 3525
$ws.Range("I134").Value = 1991.421
$ws.Range("J134").Value = 5953.1665
$ws.Range("K134").Value = 5974.263
$ws.Range("L134").Value = 17859.4995
$ws.Range("M134").Value = -3439.263
$ws.Range("N134").Value = -22929.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1294.5294
$ws.Range("I5").Value = 800.13336
$ws.Range("J5").Value = 5002.5
$ws.Range("K5").Value = 2400.40008
$ws.Range("L5").Value = 15007.5
$ws.Range("M5").Value = -2288.40008
$ws.Range("N5").Value = -15231.5

$ws.Range("H11").Value = 216111.42
$ws.Range("I11").Value = 300456
$ws.Range("K11").Value = 901368
$ws.Range("M11").Value = -901228

$ws.Range("H18").Value = 676.0909
$ws.Range("I18").Value = 459.85715
$ws.Range("J18").Value = 1054.5
$ws.Range("K18").Value = 1379.57145
$ws.Range("L18").Value = 3163.5
$ws.Range("M18").Value = -1210.57145
$ws.Range("N18").Value = -3501.5

$ws.Range("H32").Value = 99.8
$ws.Range("J32").Value = 100
$ws.Range("L32").Value = 300
$ws.Range("N32").Value = -866

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H86").Value = 662.5
$ws.Range("I86").Value = 261.5
$ws.Range("K86").Value = 784.5
$ws.Range("M86").Value = 401.5

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H89").Value = 662.5
$ws.Range("I89").Value = 261.5
$ws.Range("K89").Value = 2353.5
$ws.Range("M89").Value = 3574.5

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H113").Value = 5007.381
$ws.Range("I113").Value = 1075.4
$ws.Range("J113").Value = 6236.125
$ws.Range("K113").Value = 3226.2
$ws.Range("L113").Value = 18708.375
$ws.Range("M113").Value = -1056.2
$ws.Range("N113").Value = -23048.375

$ws.Range("H135").Value = 1294.5294
$ws.Range("I135").Value = 800.13336
$ws.Range("J135").Value = 5002.5
$ws.Range("K135").Value = 7201.20024
$ws.Range("L135").Value = 45022.5
$ws.Range("M135").Value = -4666.20024
$ws.Range("N135").Value = -50092.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6315.7144
$ws.Range("I102").Value = 6199.4165
$ws.Range("K102").Value = 6199.4165
$ws.Range("M102").Value = -4577.4165

$ws.Range("H113").Value = 6278.7676
$ws.Range("I113").Value = 5218.3125
$ws.Range("K113").Value = 5218.3125
$ws.Range("M113").Value = -3048.3125

$ws.Range("H122").Value = 8074462.5
$ws.Range("I122").Value = 18159364
$ws.Range("J122").Value = 6541.2
$ws.Range("K122").Value = 54478092
$ws.Range("L122").Value = 19623.6
$ws.Range("M122").Value = -54475642
$ws.Range("N122").Value = -24523.6

$ws.Range("H132").Value = 2900
$ws.Range("I132").Value = 1866
$ws.Range("J132").Value = 5197.778
$ws.Range("K132").Value = 5598
$ws.Range("L132").Value = 15593.334
$ws.Range("M132").Value = -3068
$ws.Range("N132").Value = -20653.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 63334
$ws.Range("I2").Value = 45001
$ws.Range("K2").Value = 45001
$ws.Range("M2").Value = -44889

$ws.Range("H40").Value = 3519.8823
$ws.Range("I40").Value = 3519.8823
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3519.8823
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3383.8823
$ws.Range("N40").ClearContents()

$ws.Range("H107").Value = 4574.5
$ws.Range("I107").Value = 4574.5
$ws.Range("K107").Value = 4574.5
$ws.Range("M107").Value = -2654.5

$ws.Range("H132").Value = 12504700
$ws.Range("I132").Value = 27780168
$ws.Range("J132").Value = 6590.909
$ws.Range("K132").Value = 83340504
$ws.Range("L132").Value = 19772.727
$ws.Range("M132").Value = -83337974
$ws.Range("N132").Value = -24832.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 57388.332
$ws.Range("I62").Value = 81082.336
$ws.Range("K62").Value = 81082.336
$ws.Range("M62").Value = -80458.336

$ws.Range("H65").Value = 57388.332
$ws.Range("I65").Value = 81082.336
$ws.Range("K65").Value = 405411.68
$ws.Range("M65").Value = -402291.68

$ws.Range("H113").Value = 1383.5714
$ws.Range("I113").Value = 1210.6666
$ws.Range("J113").Value = 1583.0769
$ws.Range("K113").Value = 3631.9998
$ws.Range("L113").Value = 4749.2307
$ws.Range("M113").Value = -1461.9998
$ws.Range("N113").Value = -9089.2307

$ws.Range("H124").Value = 44510.5
$ws.Range("J124").Value = 44510.5
$ws.Range("L124").Value = 44510.5
$ws.Range("N124").Value = -54330.5

$ws.Range("H136").Value = 23492548
$ws.Range("I136").Value = 41668188
$ws.Range("J136").Value = 533842.3
$ws.Range("K136").Value = 125004564
$ws.Range("L136").Value = 1601526.9
$ws.Range("M136").Value = -125002014
$ws.Range("N136").Value = -1606626.9
